$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    'Juvenal Bashirian V',
    'Ramon Borer',
    'Donna Haley',
    'Jerel Ledner',
    'Jesse Boehm',
    'Concepcion Tremblay',
    'Green Bosco',
    'Ignacio Dickens',
    'Lonie McKenzie',
    'Mrs. Tessie Stanton',
    'Raven Sipes',
    'Vladimir Schaden',
    'Amelie Renner',
    'Mr. Chauncey Flatley Jr.',
    'Rogelio Stokes',
    'Mr. Giovani Lemke III',
    'Zoie Batz',
    'Ms. Elinor Bartell',
    'Una Erdman V',
    'Fletcher Parisian',
    'Carey Sawayn',
    'Deontae Lynch',
    'Dallas Lemke Sr.',
    'Cali Bartoletti',
    'Shyanne Hagenes',
    'Lyda Welch',
    'Bradley Purdy',
    'Prof. Dennis Shields IV',
    'Thomas Romaguera',
    'Greyson Runte',
    'Ladarius Rippin',
    'Miss Kayla Stark',
    'Bridgette Halvorson',
    'Elwyn Will',
    'Leilani Luettgen',
    'Dr. Fiona Hackett',
    'Mr. Saige Buckridge',
    'Alexandrea Schmidt',
    'Dusty Cruickshank',
    'Joannie Jenkins',
    'Prof. Maximillian Wintheiser Sr.',
    'Ronny West',
    'June Emard',
    'Harley Ortiz',
    'Daphne Runolfsdottir',
    'Penelope Schneider',
    'Deion Parker MD',
    'Fannie Waelchi Jr.',
    'Ms. Rosalind Schiller',
    'Cheyanne Schaefer',
    'Guido Olson II',
    'Ms. Rosalind Vandervort V',
    'Christop McLaughlin',
    'Nikki Walker',
    'Carmine O''Conner',
    'Kendra Larson',
    'Maximillian Stiedemann',
    'Maverick Weimann',
    'Effie Klein',
    'Einar Huels',
    'Adela Murazik',
    'Osborne Ferry',
    'Burnice Marquardt',
    'Yolanda Legros II',
    'Christian Herzog',
    'Mr. Darius Mills',
    'Jordy Davis',
    'Gabriel Paucek V',
    'Jordon Thompson',
    'Ova Volkman',
    'Florida Abernathy',
    'Delmer Jakubowski IV',
    'Jacques Mitchell',
    'Mr. Geovanny Howe',
    'Marjorie Denesik I',
    'Hanna Ledner Sr.',
    'Miss Frieda Walker I',
    'Dr. Alena Borer',
    'Gavin Flatley',
    'Zula McDermott',
    'Rosamond Gaylord',
    'Ms. Shyanne Bailey Jr.',
    'Prof. Abbigail Green MD',
    'Terrill Marvin',
    'Prof. Carmine Prohaska',
    'Lelah Pfannerstill',
    'Dr. Titus Toy',
    'Jayde O''Connell',
    'Vallie Goodwin',
    'Dr. Amparo Gaylord DVM',
    'Ms. Josefa Pouros',
    'Russ Keebler',
    'Prof. Malachi Russel Sr.',
    'Dr. Zion Simonis Sr.',
    'Beth Farrell',
    'Sanford Medhurst DDS',
    'Mrs. Alicia Armstrong MD',
    'Mckenna Gulgowski',
    'Mollie Ondricka',
    'Scot Will I'
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $names[$i]
}

